$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.218.91"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.830.24"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'243.07"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'0.6161"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07336"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2917"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "'23.24"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "'0.07660"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.845.56"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'4.985"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'0.6727"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'82.52"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "'0.000008933"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "29.212.95"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "2.090.16"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'236.85"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'7.390"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'158.76"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "'8.535"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "'0.1389"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'17.64"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").Value = "'1.495"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'0.05791"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "'1.235"
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("D32").Value = "'4.085"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'4.102"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'1.854"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.7202"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "'2.614"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'2.864"
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("D39").Value = "1.221.41"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'0.01763"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.203"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9069"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "2.008.58"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").Value = "'101.97"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "'0.5049"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000119"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1178"
$ws.Range("E49").Value = "  +7.02%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.155"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4033"
$ws.Range("E51").Value = "  -0.17%  "

# Cells that look like plain numbers (e.g. "0.9999") get an auto quote-prefix
# style from the assignment above; restore them to the workbook default style
# (same as untouched cell A1) so no stray formatting is introduced.
$ws.Range("D4").Style = $ws.Range("A1").Style
$ws.Range("D5").Style = $ws.Range("A1").Style
$ws.Range("D6").Style = $ws.Range("A1").Style
$ws.Range("D7").Style = $ws.Range("A1").Style
$ws.Range("D8").Style = $ws.Range("A1").Style
$ws.Range("D9").Style = $ws.Range("A1").Style
$ws.Range("D10").Style = $ws.Range("A1").Style
$ws.Range("D11").Style = $ws.Range("A1").Style
$ws.Range("D13").Style = $ws.Range("A1").Style
$ws.Range("D14").Style = $ws.Range("A1").Style
$ws.Range("D15").Style = $ws.Range("A1").Style
$ws.Range("D16").Style = $ws.Range("A1").Style
$ws.Range("D20").Style = $ws.Range("A1").Style
$ws.Range("D22").Style = $ws.Range("A1").Style
$ws.Range("D23").Style = $ws.Range("A1").Style
$ws.Range("D25").Style = $ws.Range("A1").Style
$ws.Range("D26").Style = $ws.Range("A1").Style
$ws.Range("D27").Style = $ws.Range("A1").Style
$ws.Range("D28").Style = $ws.Range("A1").Style
$ws.Range("D29").Style = $ws.Range("A1").Style
$ws.Range("D30").Style = $ws.Range("A1").Style
$ws.Range("D31").Style = $ws.Range("A1").Style
$ws.Range("D32").Style = $ws.Range("A1").Style
$ws.Range("D33").Style = $ws.Range("A1").Style
$ws.Range("D34").Style = $ws.Range("A1").Style
$ws.Range("D35").Style = $ws.Range("A1").Style
$ws.Range("D36").Style = $ws.Range("A1").Style
$ws.Range("D37").Style = $ws.Range("A1").Style
$ws.Range("D38").Style = $ws.Range("A1").Style
$ws.Range("D40").Style = $ws.Range("A1").Style
$ws.Range("D41").Style = $ws.Range("A1").Style
$ws.Range("D42").Style = $ws.Range("A1").Style
$ws.Range("D45").Style = $ws.Range("A1").Style
$ws.Range("D47").Style = $ws.Range("A1").Style
$ws.Range("D48").Style = $ws.Range("A1").Style
$ws.Range("D49").Style = $ws.Range("A1").Style
$ws.Range("D50").Style = $ws.Range("A1").Style
$ws.Range("D51").Style = $ws.Range("A1").Style
